$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update evaluation metrics (based on last 30 days data)
$ws.Range("B2").Value = 0.005080933345529242
$ws.Range("C2").Value = 0.9057459357033854
$ws.Range("D2").Value = 0.8129863656775127
$ws.Range("E2").Value = 0.00002581588366171098
$ws.Range("F2").Value = 21
$ws.Range("G2").Value = 1749.95412525557
$ws.Range("H2").Value = 41.83245301504049
$ws.Range("I2").Value = 21.63134065291425
$ws.Range("J2").Value = 0.7638187495916015

# Jumlah_Data_Evaluasi (count of evaluation rows) should be shown as a plain integer
$ws.Range("F2").NumberFormat = "0"
